$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.540.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.625.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.56%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '196.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.209'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.646'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.201.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '600.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.645.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.630.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.60%  '
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('E26').Value = '  -3.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('E30').Value = '  +8.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.117'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.922.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '534.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0463'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('E51').Value = '  +1.50%  '
